$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71: new daily entry appended to the log.
# Column A holds a date-like string ("2025/10/07") that must stay as literal
# text (matching the rest of the column), not get auto-converted into a
# date serial number by Excel's input parsing. Force a Text number format
# before assigning, then reset the style back to Normal/General so no new
# cell style is introduced (the source file's data rows carry no explicit
# style).
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "2025/10/07"
$ws.Range("A71").Style = "Normal"

$ws.Range("B71").Value = "火"
$ws.Range("C71").Value = 1
$ws.Range("D71").Value = 63
